# Add "Wins", "Losses", "Ties" season-record columns (AD, AE, AF) to Sheet1,
# one value per player row, matching each team's season record.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row (row 1): new headers, styled like the rest of the header row.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-53: same season record (90-72-0) for every player row.
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 90
    $ws.Cells.Item($r, 31).Value = 72
    $ws.Cells.Item($r, 32).Value = 0
}
